$d = $word.ActiveDocument

# 1. "I'm a CSS expert..." paragraph: add "SVG, " before "SASS, " in the
#    "...extensively with SASS, LESS..." sentence.
$d.Content.Find.Execute("extensively with SASS,", $true, $false, $false, $false, $false, $true, 1, $false, "extensively with SVG, SASS,", 2)

# 2. "Technical expertise" bullet list: change
#    "HTML5 and CSS, including LESS, SASS and Bootstrap"
#    to
#    "HTML5, CSS (including LESS, SASS and Bootstrap) and SVG"
$d.Content.Find.Execute("HTML5 and CSS, including LESS, SASS and Bootstrap", $true, $false, $false, $false, $false, $true, 1, $false, "HTML5, CSS (including LESS, SASS and Bootstrap) and SVG", 2)
